# "dashboard and course updated"
#
# Updates the course-tracker table on Sheet1:
#  - Adds ATOC Week 5 / Week 6 progress (column G)
#  - Adds AOS Week 3 (re-entered), Week 4, Week 5 progress (column H)
#  - Adds STQA Week 2 / Week 3 progress (column I)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column G : ATOC -------------------------------------------------
$gWeek5 = @"
Topic: name- Week 5, lectures- 4 Lectures, duration- 01:23;
Video: link- https://drive.google.com/file/d/1MFTDv0iYQzWbW_4VCX7PgDr3-UYf4800/preview, name- ATOC Week#5Part1.mp4, duration- 00:27;
Video: link- https://drive.google.com/file/d/17RximV5JzVFI3eu-9ir6eFv8s_kbdfAe/preview, name- ATOC Week#5Part2.mp4, duration- 00:01;
Video: link- https://drive.google.com/file/d/1pSGNEYEJkci8ely09oc-rdTtnI2ljeqE/preview, name- ATOC Week#5Part3.mp4, duration- 00:53;
Video: link- https://drive.google.com/file/d/1NJH2YpnIsmDwipOJSgILBo_Ta0u6klH4/preview, name- ATOC Week#5Part4.mp4, duration- 00:02;
Note: heading- Presentation in next class;
"@

$gWeek6 = @"
Topic: name- Week 6, lectures- 3 Lectures, duration- 01:34;
Video: link- https://drive.google.com/file/d/1XSMfs_FmJKSzxEFjuYD_4Wu-qYSvFrai/preview, name- ATOC Week#6Part1.mp4, duration- 00:18;
Video: link- https://drive.google.com/file/d/1ukdXyaYhnJEsV2N4Y8-I6metWzm4_qqu/preview, name- ATOC Week#6Part2.mp4, duration- 01:14;
Video: link- https://drive.google.com/file/d/1P6PLNm_nXYdaUmOX-W4U0IwYRNgbp8FR/preview, name- ATOC Week#6Part3.mp4, duration- 00:02;
"@

$ws.Range("G7").Value = $gWeek5
$ws.Range("G7").WrapText = $true
$ws.Range("G7").VerticalAlignment = -4108
$ws.Range("G7").HorizontalAlignment = -4131
$ws.Range("G7").Font.Bold = $true

$ws.Range("G8").Value = $gWeek6
$ws.Range("G8").WrapText = $true
$ws.Range("G8").VerticalAlignment = -4108
$ws.Range("G8").HorizontalAlignment = -4131
$ws.Range("G8").Font.Bold = $true

# --- Column I : STQA ---------------------------------------------------
$iWeek2 = @"
Topic: name- Week 2, lectures- 1 Lecture, duration- 00:35;
Video: link- https://drive.google.com/file/d/1Q6d8GiRyOwbBX8aSjLix4pxttUs9xB_-/preview, name- STQA Week#2part1.mp4, duration- 00:35;
"@

$iWeek3 = @"
Topic: name- Week 3, lectures- 1 Lecture, duration- 01:00;
Video: link- https://drive.google.com/file/d/1ULwy1htskQHkuXh684Rr94S-TJnuLMWa/preview, name- STQA Week#3part1.mp4, duration- 01:00;
"@

$ws.Range("I4").Value = $iWeek2
$ws.Range("I4").WrapText = $true
$ws.Range("I4").VerticalAlignment = -4108
$ws.Range("I4").HorizontalAlignment = -4131
$ws.Range("I4").Font.Bold = $true

$ws.Range("I5").Value = $iWeek3
$ws.Range("I5").WrapText = $true
$ws.Range("I5").VerticalAlignment = -4108
$ws.Range("I5").HorizontalAlignment = -4131
$ws.Range("I5").Font.Bold = $true

# --- Column H : AOS ------------------------------------------------
$hWeek3 = @"
Topic: name- Week 3, lectures- 2 Lectures, duration- 01:38;
Video: link- https://drive.google.com/file/d/1CEA0dNMZFnAwwcInBt37frFSakNJSyuO/preview, name- AOS Week#3part1.mp4, duration- 01:06;
Video: link- https://drive.google.com/file/d/19aBooNe8gvSiXqm-yxEMN2AHg0KthhYP/preview, name- AOS Week#3part2.mp4, duration- 00:32;
"@

$hWeek4 = @"
Topic: name- Week 4, lectures- 3 Lectures, duration- 01:23;
Video: link- https://drive.google.com/file/d/1-0mGSVBAilQy22W6LntqSRnrUnam_9OS/preview, name- AOS Week4part1.mp4, duration- 00:09;
Video: link- https://drive.google.com/file/d/1eVs08UR6fAFRWO1dj7ejMPo19JCGVGXB/preview, name- AOS Week#4part2.mp4, duration- 01:08;
Video: link- https://drive.google.com/file/d/1a99Py9JBVIGX8qyeg-GS36N1BRAxhTbl/preview, name- AOS Week#4part3.mp4, duration- 00:06;
"@

$hWeek5 = @"
Topic: name- Week 5, lectures- 2 Lectures, duration- 01:09;
Video: link- https://drive.google.com/file/d/1p65AFQnFqjgadbb6Z6iOBRrxk4hrAPYU/preview, name- AOS Week5part1.mp4, duration- 01:01;
Video: link- https://drive.google.com/file/d/1v73JQCmN0GASJYm3syH_MdRRcuVtO8kM/preview, name- AOS Week#5part2.mp4, duration- 00:08;
"@

$ws.Range("H5").Value = $hWeek3
$ws.Range("H5").WrapText = $true
$ws.Range("H5").VerticalAlignment = -4107
$ws.Range("H5").HorizontalAlignment = -4131
$ws.Range("H5").Font.Bold = $true

$ws.Range("H6").Value = $hWeek4
$ws.Range("H6").WrapText = $true
$ws.Range("H6").VerticalAlignment = -4107
$ws.Range("H6").HorizontalAlignment = -4131
$ws.Range("H6").Font.Bold = $true

$ws.Range("H7").Value = $hWeek5
$ws.Range("H7").WrapText = $true
$ws.Range("H7").VerticalAlignment = -4107
$ws.Range("H7").HorizontalAlignment = -4131
$ws.Range("H7").Font.Bold = $true

# --- Update the active selection to reflect where the author left off ---
$ws.Range("H5").Select()
